# Updated cryptos list on Wed Nov 29 18:00:42 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row.
# Cells whose new text looks like a plain number (e.g. "227.68") get their
# NumberFormat forced to "@" (Text) first so Excel stores the literal
# string instead of silently coercing it to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.854.40"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.034.89"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.68"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.26"
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0816"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "2.337.43"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.16"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.755"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "2.052.48"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "37.864.13"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.05"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.37"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.96"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.90"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("E30").Value = "  -6.88%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.42"
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "1.535.67"
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.92"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.78"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("E47").Value = "  -3.41%  "
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "2.227.44"
$ws.Range("E51").Value = "  -1.66%  "
